$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Raum"
$ws.Range("B1").Value = "Kapazität"
$ws.Range("B2").Value = 15
$ws.Range("B10").Value = 20
[void]$ws.Range("B15").Select()
